$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D header ("Coef"), matching the bold style used by A1:C1
$ws.Range("D1").Value = "Coef"
$ws.Range("D1").Font.Bold = $true

# Fill D11:D37 with 1 (coefficient column added alongside existing data rows)
$ws.Range("D11:D37").Value = 1

# Update the view: scrolled down so row 19 is the top visible row, with
# D35 left as the active/selected cell
$excel.Goto($ws.Range("A19"), $true)
$ws.Range("D35").Select()
